# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from the last
# existing header cell (G1) onto the new header cell (H1) so the new
# column matches the look of the rest of the header row.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set the new header label and the row-2 value for the new column.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
